# Connected mongoDB & fixed tag loading late and not being scraped issue
#
# - A2 gets the updated (tracked) Google Project Management link
# - A3 gets a newly-scraped course link (was blank before)
# - A4 loses its leftover "Hyperlink" cell style (copy A3's plain style over it)
# - The now-unused "Hyperlink" named cell style is removed from the workbook
# - Final selection lands on A5 (last cell touched by the scraper)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://www.coursera.org/professional-certificates/google-project-management?irclickid=XFtxnQ1ldxyPROYzPS3u9VCcUkFQG639OQKVxM0&irgwc=1&utm_medium=partners&utm_source=impact&utm_campaign=4807567&utm_content=b2c"
$ws.Range("A3").Value = "https://www.coursera.org/learn/meem-agile-project-management?irclickid=XFtxnQ1ldxyPROYzPS3u9VCcUkFQG6RVOQKVxM0&irgwc=1&utm_medium=partners&utm_source=impact&utm_campaign=4807567&utm_content=b2c"

# A4 previously carried the built-in "Hyperlink" cell style (underlined,
# theme-colored font) left over from a removed hyperlink; bring it back to
# the plain style used by its neighbors.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "Hyperlink" cell style is no longer referenced by any cell - delete it
# from the workbook's style gallery.
$wb.Styles.Item("Hyperlink").Delete()

$ws.Range("A5").Select()
